# STM32G474 pinout.xlsx - add a "UART" overview sheet between "Timer" and "ADC",
# and note the pin change for L26 on the "ADC" sheet (cap3k / anti-aliasing edit).

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "UART" worksheet right before the current active sheet
#        ("ADC"), matching where Excel drops a freshly-inserted sheet. -----------
$uart = $wb.Worksheets.Add()
$uart.Name = "UART"

# --- 2. Fill in the UART pin-assignment overview table ---------------------------
$uart.Range("A1").Value = "UART Module"
$uart.Range("B1").Value = "Available RX"
$uart.Range("C1").Value = "Available TX"

$uart.Range("A2").Value = "USART1"
$uart.Range("B2").Value = "PC5, PE1, (PB7)"
$uart.Range("C2").Value = "PC4, PE0, (PB6)"

$uart.Range("A3").Value = "USART2"
$uart.Range("B3").Value = "PD6, (PB4)"
$uart.Range("C3").Value = "PD5"

$uart.Range("A4").Value = "UART3"
$uart.Range("B4").Value = "PB8"
$uart.Range("C4").Value = "PB9"

$uart.Range("A5").Value = "UART4"
$uart.Range("B5").Value = "Blocked by SPI"
$uart.Range("C5").Value = "Blocked by SPI"

$uart.Range("A6").Value = "UART5"
$uart.Range("B6").Value = "PD2"
$uart.Range("C6").Value = "Blocked by SPI"

# Row 3 (USART2, the module already routed/used) is highlighted green.
$uart.Range("A3:C3").Interior.Color = 5296274   # RGB(146,208,80) -> 0x92D050 (BGR packed)

# Column widths roughly matching the "best fit" widths Excel computed.
$uart.Columns.Item(1).ColumnWidth = 11.666666666666666
$uart.Columns.Item(2).ColumnWidth = 12.666666666666666
$uart.Columns.Item(3).ColumnWidth = 12.666666666666666

# Leave the cursor where the author left it on this sheet.
$uart.Range("C12").Select()

# --- 3. Re-activate "ADC" so it stays the active/visible tab (only its position
#        in the tab strip moved because of the newly inserted sheet). -----------
$adc = $wb.Worksheets.Item("ADC")
$adc.Activate()

# Pin "I2C3_SDA"/PB8 row: routing moved the chosen alternative from option 1 to 5,
# and a remark was added to flag it changed during layouting.
$adc.Range("L26").Value = 5
$adc.Range("M26").Value = "Pin changed during layouting"

$adc.Range("M27").Select()
